# Update the "dSF" column (column F) values for rows 2-40 on the active sheet,
# per repulled data / recalculated means.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 10
    3  = -1
    4  = -4
    5  = 1
    6  = -1
    7  = 1
    9  = 6
    10 = 2
    11 = 3
    12 = 1
    13 = -2
    14 = -1
    15 = -1
    16 = 2
    17 = 1
    18 = -4
    20 = 2
    21 = 2
    22 = -1
    23 = 2
    25 = 2
    26 = -1
    27 = -4
    28 = -3
    29 = 5
    30 = 1
    31 = 9
    32 = -1
    33 = -3
    34 = -1
    35 = -6
    36 = -4
    37 = -1
    38 = -1
    39 = -1
    40 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
